$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in row 2 and row 3 for the columns that actually changed
# (A, B, E, F, G, H, I, Q, R) - these two records got their row order swapped.
$numericCols = @("A", "B", "E", "Q", "R")
$textCols = @("F", "G", "H")

foreach ($col in $numericCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

foreach ($col in $textCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# Column I holds numeric-looking values but is stored as text in the
# workbook, so force the cell format to Text before writing so Excel
# doesn't silently re-type it as a number.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I2").Value = "$i3"
$ws.Range("I3").Value = "$i2"
